# PerformanceTestData.xlsx correction
# "fix: heavy loadtest data corrected according to confluence Performance Test Data.xls"
#
# Corrects several input cells on the "heavy" sheet (ConnectorLight / row 6 is
# effectively zeroed-out, ConnectorMedium / row 7's ratio moves from 24% to
# 25%), lets the workbook's formulas recalculate the dependent totals, and
# updates the sheet selections/active sheet to match where the author ended
# up after editing.

$wb = $excel.ActiveWorkbook

$heavy = $wb.Worksheets.Item("heavy")
$test  = $wb.Worksheets.Item("test")

# --- Row 6 (ConnectorLight) : zero everything out -------------------------
$heavy.Range("E6").Value = 0
$heavy.Range("G6").Value = 0
# H6 used to be the formula "=100+K6"; it's now a plain literal 0.
$heavy.Range("H6").Value = 0
$heavy.Range("M6").Value = 0

# --- Row 7 (ConnectorMedium) : ratio 24% -> 25% ----------------------------
$heavy.Range("G7").Value = 0.25
$heavy.Range("M7").Value = 0.25

# All the dependent formula cells (R2, S2, U2, X2, I5, F6, I6, J6, F7, I7, ...)
# recalculate automatically from the inputs above.

# --- Selections / active sheet --------------------------------------------
# Author ended up on "heavy", scrolled to column D, with D8 selected.
$heavy.Activate()
$heavy.Range("D8").Select()

# "test" sheet is no longer the active tab; its lingering selection moved to E6.
$test.Range("E6").Select()

# Leave "heavy" as the active sheet/tab.
$heavy.Activate()
